$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update row 1 header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 (CON) values (B2:E2)
$ws.Range("B2").Value = 7.3099930753124456
$ws.Range("C2").Value = 5.0637421277820867
$ws.Range("D2").Value = 6.1868900718345197
$ws.Range("E2").Value = 6.9414908877550401

# Update row 3 (STR) values (B3:E3)
$ws.Range("B3").Value = 6.2880561177984298
$ws.Range("C3").Value = 4.4550762181419969
$ws.Range("D3").Value = 6.0647960868018229
$ws.Range("E3").Value = 7.2646165724020548

# Update selection to reflect the new active range B1:E3
$ws.Range("B1:E3").Select()
